$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh. Column D holds prices that are stored as
# plain text (so values like "3.90" or "0.0000143" keep their exact
# formatting instead of being normalized as numbers). Force text entry via
# a temporary "@" number format, then restore the original style so the
# cell format is left untouched, matching the source data.

$c = $ws.Range("D2")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "61.853.79"
$c.Style = $origStyle
$ws.Range("E2").Value = "  -2.30%  "
$c = $ws.Range("D3")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.576.19"
$c.Style = $origStyle
$ws.Range("E3").Value = "  -3.93%  "
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "549.26"
$c.Style = $origStyle
$ws.Range("E5").Value = "  -0.74%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "153.84"
$c.Style = $origStyle
$ws.Range("E6").Value = "  -2.66%  "
$c = $ws.Range("D7")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle
$ws.Range("E7").Value = "  -0.07%  "
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.591"
$c.Style = $origStyle
$ws.Range("E8").Value = "  +0.74%  "
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.103"
$c.Style = $origStyle
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -0.52%  "
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.49"
$c.Style = $origStyle
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("E12").Value = "  -0.82%  "
$c = $ws.Range("D13")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.032.66"
$c.Style = $origStyle
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("E14").Value = "  -3.55%  "
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "61.761.20"
$c.Style = $origStyle
$ws.Range("E15").Value = "  -2.24%  "
$c = $ws.Range("D16")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0000143"
$c.Style = $origStyle
$ws.Range("E16").Value = "  -1.18%  "
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.578.22"
$c.Style = $origStyle
$ws.Range("E17").Value = "  -3.99%  "
$c = $ws.Range("D18")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "11.54"
$c.Style = $origStyle
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("E19").Value = "  -0.96%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "337.78"
$c.Style = $origStyle
$ws.Range("E20").Value = "  -2.29%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.99"
$c.Style = $origStyle
$ws.Range("E21").Value = "  -4.84%  "
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle
$ws.Range("E22").Value = "  +0.20%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.488"
$c.Style = $origStyle
$ws.Range("E23").Value = "  -3.55%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "63.42"
$c.Style = $origStyle
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  -0.93%  "
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle
$ws.Range("E26").Value = "  -0.09%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.09"
$c.Style = $origStyle
$ws.Range("E27").Value = "  -0.79%  "
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0₃0828"
$c.Style = $origStyle
$ws.Range("E28").Value = "  -3.09%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.13"
$c.Style = $origStyle
$ws.Range("E29").Value = "  +1.84%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("E31").Value = "  -2.62%  "
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "161.86"
$c.Style = $origStyle
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("E33").Value = "  +0.01%  "
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.82"
$c.Style = $origStyle
$ws.Range("E34").Value = "  +0.06%  "
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "19.08"
$c.Style = $origStyle
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  -2.19%  "
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.77"
$c.Style = $origStyle
$ws.Range("E37").Value = "  -0.06%  "
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.01"
$c.Style = $origStyle
$ws.Range("E38").Value = "  -1.74%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "323.66"
$c.Style = $origStyle
$ws.Range("E39").Value = "  -4.87%  "
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.892"
$c.Style = $origStyle
$ws.Range("E40").Value = "  -4.79%  "
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.90"
$c.Style = $origStyle
$ws.Range("E41").Value = "  -0.64%  "
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "37.44"
$c.Style = $origStyle
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("E43").Value = "  -0.96%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0962"
$c.Style = $origStyle
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "19.42"
$c.Style = $origStyle
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0541"
$c.Style = $origStyle
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("E50").Value = "  -1.40%  "
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.036.67"
$c.Style = $origStyle
$ws.Range("E51").Value = "  -2.67%  "
